# Insert a new data row at row 258 (pushing the existing rows 258:279 down
# to 259:280) and populate it with the new weekly price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 258:279 down one row to make room for the new record.
$ws.Rows("258").Insert()

# Populate the newly inserted row 258 with the new weekly record.
$ws.Range("A258").Value = 5
$ws.Range("B258").Value = "Macroferia Regional de Talca"
$ws.Range("C258").Value = "Maule"
$ws.Range("D258").Value = 44931
$ws.Range("E258").Value = 7
$ws.Range("F258").Value = 100112024
$ws.Range("G258").Value = "Choclo"
$ws.Range("H258").Value = "Choclero"
$ws.Range("I258").Value = "Primera"
$ws.Range("J258").Value = 40000
$ws.Range("K258").Value = 230
$ws.Range("L258").Value = 250
$ws.Range("M258").Value = 240
$ws.Range("N258").Value = "`$/unidad"
$ws.Range("O258").Value = "Región del Maule"
$ws.Range("P258").Value = 240
$ws.Range("Q258").Value = 1
$ws.Range("R258").Value = "Hortaliza"
